# "Add files via upload" - MLB Stats 2018.xlsx
# Rename the batting average header from "BA" to "AVG", highlight several
# key columns on the "2018 League Hitting" sheet, and add a new "wRC+"
# column with per-team data and a league-average formula.

$wb = $excel.ActiveWorkbook

$pitching = $wb.Worksheets.Item("2018 League Pitching")
$hitting  = $wb.Worksheets.Item("2018 League Hitting")

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------
# 1. Rename "BA" -> "AVG" on the hitting sheet header row.
# ---------------------------------------------------------------------
$hitting.Range("R1").Value = "AVG"

# ---------------------------------------------------------------------
# 2. Highlight (yellow fill) a handful of header cells: G, PA, SB, AVG,
#    OBP, SLG, OPS+. Re-use the existing highlighted-header format from
#    the pitching sheet (column H1) so the same style slot is reused.
# ---------------------------------------------------------------------
$pitching.Range("H1").Copy() | Out-Null
$highlightCells = @("E1", "F1", "N1", "R1", "S1", "T1", "V1")
foreach ($addr in $highlightCells) {
    $hitting.Range($addr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# 3. Add the new "wRC+" column (AE) after WAR (AD).
# ---------------------------------------------------------------------
$hitting.Range("AD1").Copy() | Out-Null
$hitting.Range("AE1").PasteSpecial(-4122) | Out-Null
$hitting.Range("AE1").Value = "wRC+"

$wrcValues = @(88, 96, 87, 111, 100, 92, 94, 106, 88, 83, 110, 87, 98, 111, 83, 100, 96, 94, 112, 112, 91, 95, 84, 100, 82, 99, 107, 89, 100, 101)

$pitching.Range("AM3").Copy() | Out-Null
for ($i = 0; $i -lt $wrcValues.Length; $i++) {
    $row = 2 + $i
    $cell = $hitting.Range("AE" + $row)
    $cell.PasteSpecial(-4122) | Out-Null
    $cell.Value = $wrcValues[$i]
}

$hitting.Range("AD32").Copy() | Out-Null
$hitting.Range("AE32").PasteSpecial(-4122) | Out-Null
$hitting.Range("AE32").Formula = "=AVERAGE(AE2:AE31)"

# ---------------------------------------------------------------------
# 4. Touch AF1/AG1 (empty, but formatted) so the sheet's used range
#    extends out to column AG, matching the author's selection drag.
# ---------------------------------------------------------------------
$hitting.Range("AD1").Copy() | Out-Null
$hitting.Range("AF1").PasteSpecial(-4122) | Out-Null
$hitting.Range("AG1").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 5. Update the hitting sheet view (scrolled right, AG1 selected).
# ---------------------------------------------------------------------
$hitting.Activate()
try { $hitting.Application.ActiveWindow.TopLeftCell = $hitting.Range("M1") } catch {}
$hitting.Range("AG1").Select() | Out-Null
